$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats
$xlPasteFormats = -4122

# --- Column S (2022) data, mirroring column R's styling for each row ---

$ws.Range("R4").Copy()
$ws.Range("S4").PasteSpecial($xlPasteFormats)
$ws.Range("S4").Value = 2022

$ws.Range("R5").Copy()
$ws.Range("S5").PasteSpecial($xlPasteFormats)
$ws.Range("S5").Value = 4.9000000000000004

$ws.Range("R6").Copy()
$ws.Range("S6").PasteSpecial($xlPasteFormats)
$ws.Range("S6").Value = 6.1

$ws.Range("R7").Copy()
$ws.Range("S7").PasteSpecial($xlPasteFormats)
$ws.Range("S7").Value = 4

# Row 8 is a blank "header" row (like D8:Q8); it gets a new bold+italic
# sz9 Times New Roman style instead of copying R8's plain style.
$ws.Range("R6").Copy()
$ws.Range("S8").PasteSpecial($xlPasteFormats)
$ws.Range("S8").Font.Bold = $true
$ws.Range("S8").Font.Italic = $true

$ws.Range("R9").Copy()
$ws.Range("S9").PasteSpecial($xlPasteFormats)
$ws.Range("S9").Value = 6.1

$ws.Range("R10").Copy()
$ws.Range("S10").PasteSpecial($xlPasteFormats)
$ws.Range("S10").Value = 12.4

$ws.Range("R11").Copy()
$ws.Range("S11").PasteSpecial($xlPasteFormats)
$ws.Range("S11").Value = 3.2

$ws.Range("R12").Copy()
$ws.Range("S12").PasteSpecial($xlPasteFormats)
$ws.Range("S12").Value = 10.8

$ws.Range("R13").Copy()
$ws.Range("S13").PasteSpecial($xlPasteFormats)
$ws.Range("S13").Value = 14.6

$ws.Range("R14").Copy()
$ws.Range("S14").PasteSpecial($xlPasteFormats)
$ws.Range("S14").Value = 8.5

$ws.Range("R15").Copy()
$ws.Range("S15").PasteSpecial($xlPasteFormats)
$ws.Range("S15").Value = 5.5

$ws.Range("R16").Copy()
$ws.Range("S16").PasteSpecial($xlPasteFormats)
$ws.Range("S16").Value = 7.1

$ws.Range("R17").Copy()
$ws.Range("S17").PasteSpecial($xlPasteFormats)
$ws.Range("S17").Value = 4.4000000000000004

$ws.Range("R18").Copy()
$ws.Range("S18").PasteSpecial($xlPasteFormats)
$ws.Range("S18").Value = 5.8

$ws.Range("R19").Copy()
$ws.Range("S19").PasteSpecial($xlPasteFormats)
$ws.Range("S19").Value = 11.6

$ws.Range("R20").Copy()
$ws.Range("S20").PasteSpecial($xlPasteFormats)
$ws.Range("S20").Value = 3.1

$ws.Range("R21").Copy()
$ws.Range("S21").PasteSpecial($xlPasteFormats)
$ws.Range("S21").Value = 1.5

$ws.Range("R22").Copy()
$ws.Range("S22").PasteSpecial($xlPasteFormats)
$ws.Range("S22").Value = 2.2999999999999998

$ws.Range("R23").Copy()
$ws.Range("S23").PasteSpecial($xlPasteFormats)
$ws.Range("S23").Value = 1

$ws.Range("R24").Copy()
$ws.Range("S24").PasteSpecial($xlPasteFormats)
$ws.Range("S24").Value = 2.2999999999999998

$ws.Range("R25").Copy()
$ws.Range("S25").PasteSpecial($xlPasteFormats)
$ws.Range("S25").Value = 3.3

$ws.Range("R26").Copy()
$ws.Range("S26").PasteSpecial($xlPasteFormats)
$ws.Range("S26").Value = 1.6

$ws.Range("R27").Copy()
$ws.Range("S27").PasteSpecial($xlPasteFormats)
$ws.Range("S27").Value = 4.5999999999999996

$ws.Range("R28").Copy()
$ws.Range("S28").PasteSpecial($xlPasteFormats)
$ws.Range("S28").Value = 4.4000000000000004

$ws.Range("R29").Copy()
$ws.Range("S29").PasteSpecial($xlPasteFormats)
$ws.Range("S29").Value = 4.7

$ws.Range("R30").Copy()
$ws.Range("S30").PasteSpecial($xlPasteFormats)
$ws.Range("S30").Value = 4

$ws.Range("R31").Copy()
$ws.Range("S31").PasteSpecial($xlPasteFormats)
$ws.Range("S31").Value = 3.2

$ws.Range("R32").Copy()
$ws.Range("S32").PasteSpecial($xlPasteFormats)
$ws.Range("S32").Value = 4.7

$ws.Range("R33").Copy()
$ws.Range("S33").PasteSpecial($xlPasteFormats)
$ws.Range("S33").Value = 2.6

$ws.Range("R34").Copy()
$ws.Range("S34").PasteSpecial($xlPasteFormats)
$ws.Range("S34").Value = 3.3

$ws.Range("R35").Copy()
$ws.Range("S35").PasteSpecial($xlPasteFormats)
$ws.Range("S35").Value = 2.2000000000000002

# Row 36 is the other blank "header" row (like D36:Q36); same treatment as row 8.
$ws.Range("R6").Copy()
$ws.Range("S36").PasteSpecial($xlPasteFormats)
$ws.Range("S36").Font.Bold = $true
$ws.Range("S36").Font.Italic = $true

$ws.Range("R37").Copy()
$ws.Range("S37").PasteSpecial($xlPasteFormats)
$ws.Range("S37").Value = 13.2

$ws.Range("R38").Copy()
$ws.Range("S38").PasteSpecial($xlPasteFormats)
$ws.Range("S38").Value = 7.5

$ws.Range("R39").Copy()
$ws.Range("S39").PasteSpecial($xlPasteFormats)
$ws.Range("S39").Value = 4.0999999999999996

$ws.Range("R40").Copy()
$ws.Range("S40").PasteSpecial($xlPasteFormats)
$ws.Range("S40").Value = 4.3

$ws.Range("R41").Copy()
$ws.Range("S41").PasteSpecial($xlPasteFormats)
$ws.Range("S41").Value = 2.6

$ws.Range("R42").Copy()
$ws.Range("S42").PasteSpecial($xlPasteFormats)
$ws.Range("S42").Value = 1

# Row 43 footnote marker "…" (shared string, same as R43)
$ws.Range("R43").Copy()
$ws.Range("S43").PasteSpecial($xlPasteFormats)
$ws.Range("S43").Value = "…"

# Match the saved selection/active cell from the authored workbook.
$null = $ws.Range("T12").Select()
